# Weekly update: insert a new data row (row 4) for the latest observation,
# pushing all the existing daily readings down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 - this shifts rows 4..42 down to 5..43
# and keeps everything else (formatting, headers) intact.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new observation. The "template"
# columns (A,B,C,E,F,G,H,I,J,K,L,Q,T) are identical on every data row of this
# sheet, so we simply repeat them; the variable columns (D,M,N,O,P,R,S) carry
# the new values from this week's entry.
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7500
$ws.Range("Q4").Value = "$/bandeja 2 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 3750
$ws.Range("T4").Value = 2
